$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 105, shifting existing rows 105:182 down to 106:183
$ws.Rows.Item(105).Insert()

# Populate the new row 105 with the new record's data
$ws.Range("A105").Value = 10
$ws.Range("B105").Value = "Vega Modelo de Temuco"
$ws.Range("C105").Value = "La Araucanía"
$ws.Range("D105").Value = 44603
$ws.Range("E105").Value = 9
$ws.Range("F105").Value = "Fruta"
$ws.Range("G105").Value = 100103
$ws.Range("H105").Value = "Frutos de hueso (carozo)"
$ws.Range("I105").Value = 100103002
$ws.Range("J105").Value = "Ciruela"
$ws.Range("K105").Value = "Black Amber"
$ws.Range("L105").Value = "Primera"
$ws.Range("M105").Value = 160
$ws.Range("N105").Value = 13000
$ws.Range("O105").Value = 14000
$ws.Range("P105").Value = 13625
$ws.Range("Q105").Value = "$/bandeja 18 kilos granel"
$ws.Range("R105").Value = "Región de O'Higgins"
$ws.Range("S105").Value = 757
$ws.Range("T105").Value = 18
